$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new process word ------------------------------------------------
# Append "markieren" (new Prozesswort) with its definition as a new row at the
# bottom of the table, using the same formatting as the rest of the rows:
# column A = wrap text + top-aligned, column B = wrap text.
$newRow = 19
$ws.Cells.Item($newRow, 1).Value = "markieren"
$ws.Cells.Item($newRow, 2).Value = "das Selektieren einer dargebotenen Option auf der Nutzeroberfläche durch optische Hervorhebung"
$ws.Cells.Item($newRow, 1).WrapText = $true
$ws.Cells.Item($newRow, 1).VerticalAlignment = -4160
$ws.Cells.Item($newRow, 2).WrapText = $true

# --- Re-sort the table alphabetically ----------------------------------------
# The table (A2:B19) is kept sorted alphabetically by the Prozesswort column,
# so re-apply the sort to move "markieren" into its correct position.
$sortRange = $ws.Range("A2:B19")
$sortKey = $ws.Range("A2:A19")
$sortRange.Sort($sortKey)

# --- Refresh row heights (word-wrap reflow after the insert/sort) -----------
$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(4).RowHeight = 57.6
$ws.Rows.Item(5).RowHeight = 28.8
$ws.Rows.Item(9).RowHeight = 57.6
$ws.Rows.Item(11).RowHeight = 43.2
$ws.Rows.Item(12).RowHeight = 28.8
$ws.Rows.Item(16).RowHeight = 43.2
$ws.Rows.Item(17).RowHeight = 43.8
$ws.Rows.Item(18).RowHeight = 14.4
$ws.Rows.Item(19).RowHeight = 28.8

# --- Update the active selection ---------------------------------------------
$ws.Range("F17").Select()
